$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gaz")

# The sheet currently has a header row (row 1) and a single data row
# (row 2: "2025-06-17", -, -, -). We need to insert a new data row
# above it for "2025-06-16", pushing the existing row down to row 3.

# 1) Copy the existing row 2 values down into row 3 first.
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = $ws.Cells.Item(2, 1).Value2
$ws.Cells.Item(3, 2).Value = $ws.Cells.Item(2, 2).Value2
$ws.Cells.Item(3, 3).Value = $ws.Cells.Item(2, 3).Value2
$ws.Cells.Item(3, 4).Value = $ws.Cells.Item(2, 4).Value2
# Reset to the default (unstyled) cell style, matching the other data
# rows in this sheet which carry no explicit formatting.
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(3, 2).Style = "Normal"
$ws.Cells.Item(3, 3).Style = "Normal"
$ws.Cells.Item(3, 4).Style = "Normal"

# 2) Overwrite row 2 with the new date entry. Force text formatting
# first so the date string is kept as literal text (e.g.
# "2025-06-16") instead of being auto-converted to a date serial
# number, then reset the style so no explicit formatting is left on
# the cell (same as the other data rows).
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-06-16"
$ws.Cells.Item(2, 2).Value = "-"
$ws.Cells.Item(2, 3).Value = "-"
$ws.Cells.Item(2, 4).Value = "-"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 2).Style = "Normal"
$ws.Cells.Item(2, 3).Style = "Normal"
$ws.Cells.Item(2, 4).Style = "Normal"
